$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Sheet"

# Update header cell values
$ws.Range("B1").Value = "Local"
$ws.Range("C1").Value = "W"
$ws.Range("D1").Value = "H"
